# BAO CAO TUAN 4.docx - apply commit changes
#
# 1) Center the two title paragraphs and wrap them (plus the paragraph
#    mark that separates them) in the "_GoBack" bookmark - this also
#    relocates _GoBack away from its old spot (Word only keeps a single
#    bookmark per name, so re-adding it elsewhere removes the old one).
# 2) Fix up the footer's cached PAGE field result (2 -> 1).

$d = $word.ActiveDocument

# --- Center the title paragraphs -------------------------------------
$d.Paragraphs.Item(1).Alignment = 1   # wdAlignParagraphCenter
$d.Paragraphs.Item(2).Alignment = 1   # wdAlignParagraphCenter

# --- Move the "_GoBack" bookmark to span both title paragraphs -------
# Paragraph 1 starts the document; Paragraph 2 ends right before
# Heading2 ("Khao sat nhu cau nguoi dung") begins, so this span covers
# exactly the two title lines (including the paragraph mark between
# them). Adding a bookmark with a name that already exists elsewhere
# in the document moves it here and removes the prior occurrence.
$titleEnd = $d.Paragraphs.Item(2).Range.End
$titleRange = $d.Range(0, $titleEnd)
$d.Bookmarks.Add("_GoBack", $titleRange)

# --- Update the footer page-number field's cached result -------------
$footer = $d.Sections.Item(1).Footers.Item(1)
$footerRange = $footer.Range.Duplicate
$footerRange.Find.Execute("2", $true, $false, $false, $false, $false, `
                           $true, 1, $false, "1", 2)
